$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value looks like a plain number (e.g. "559.07") must be
# forced to Text format first, otherwise Excel auto-converts the assigned
# string into a numeric value -- the source sheet keeps these as literal
# text strings (e.g. "1.00", "6.00").
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '65.979.00'
$ws.Range('E2').Value = '  -4.93%  '
$ws.Range('D3').Value = '3.335.54'
$ws.Range('E3').Value = '  -5.79%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '559.07'
$ws.Range('E5').Value = '  -4.34%  '
$ws.Range('D6').Value = '181.38'
$ws.Range('E6').Value = '  -7.94%  '
$ws.Range('E7').Value = '  +0.12%  '
$ws.Range('D8').Value = '0.587'
$ws.Range('E8').Value = '  -3.64%  '
$ws.Range('D9').Value = '3.327.42'
$ws.Range('E9').Value = '  -5.71%  '
$ws.Range('D10').Value = '0.184'
$ws.Range('E10').Value = '  -9.67%  '
$ws.Range('D11').Value = '0.583'
$ws.Range('E11').Value = '  -7.19%  '
$ws.Range('D12').Value = '47.04'
$ws.Range('E12').Value = '  -9.12%  '
$ws.Range('D13').Value = '0.0000264'
$ws.Range('E13').Value = '  -7.81%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '3.868.37'
$ws.Range('E14').Value = '  -5.70%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '8.56'
$ws.Range('E15').Value = '  -7.13%  '
$ws.Range('D16').Value = '601.72'
$ws.Range('E16').Value = '  -9.48%  '
$ws.Range('D17').Value = '18.05'
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').Value = '66.026.01'
$ws.Range('E18').Value = '  -5.08%  '
$ws.Range('D19').Value = '3.348.87'
$ws.Range('E19').Value = '  -5.40%  '
$ws.Range('E20').Value = '  -3.98%  '
$ws.Range('D21').Value = '11.37'
$ws.Range('E21').Value = '  -9.05%  '
$ws.Range('D22').Value = '0.903'
$ws.Range('E22').Value = '  -6.38%  '
$ws.Range('D23').Value = '16.74'
$ws.Range('E23').Value = '  -8.94%  '
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '5.02'
$ws.Range('E24').Value = '  -5.33%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = '100.23'
$ws.Range('E25').Value = '  -4.35%  '
$ws.Range('D26').Value = '4.02'
$ws.Range('E26').Value = '  -7.97%  '
$ws.Range('D27').Value = '6.00'
$ws.Range('E27').Value = '  +0.15%  '
$ws.Range('D28').Value = '2.65'
$ws.Range('E28').Value = '  -8.52%  '
$ws.Range('D29').Value = '9.21'
$ws.Range('E29').Value = '  -9.29%  '
$ws.Range('D30').Value = '8.67'
$ws.Range('E30').Value = '  -9.55%  '
$ws.Range('D31').Value = '30.57'
$ws.Range('E31').Value = '  -7.79%  '
$ws.Range('D32').Value = '6.22'
$ws.Range('E32').Value = '  -8.11%  '
$ws.Range('D33').Value = '3.72'
$ws.Range('E33').Value = '  -15.50%  '
$ws.Range('D34').Value = '10.98'
$ws.Range('E34').Value = '  -6.55%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').Value = '0.104'
$ws.Range('E35').Value = '  -6.12%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '3.774.20'
$ws.Range('E36').Value = '  -0.46%  '
$ws.Range('B37').Value = 'OKB'
$ws.Range('C37').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D37').Value = '57.75'
$ws.Range('E37').Value = '  -6.76%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '530.22'
$ws.Range('E38').Value = '  +5.65%  '
$ws.Range('D39').Value = '0.997'
$ws.Range('E39').Value = '  -0.21%  '
$ws.Range('D40').Value = '3.39'
$ws.Range('E40').Value = '  -8.96%  '
$ws.Range('D41').Value = '0.0₃0706'
$ws.Range('E41').Value = '  -13.32%  '
$ws.Range('D42').Value = '2.64'
$ws.Range('E42').Value = '  -9.29%  '
$ws.Range('D43').Value = '0.124'
$ws.Range('E43').Value = '  -7.38%  '
$ws.Range('D44').Value = '0.339'
$ws.Range('E44').Value = '  -8.78%  '
$ws.Range('D45').Value = '31.69'
$ws.Range('E45').Value = '  -8.00%  '
$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').Value = '0.0411'
$ws.Range('E46').Value = '  -8.13%  '
$ws.Range('B47').Value = 'CoreDAO'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range('D47').Value = '3.13'
$ws.Range('E47').Value = '  +18.15%  '
$ws.Range('D48').Value = '3.21'
$ws.Range('E48').Value = '  -5.28%  '
$ws.Range('B49').Value = 'Stellar'
$ws.Range('C49').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D49').Value = '0.129'
$ws.Range('E49').Value = '  -5.53%  '
$ws.Range('B50').Value = 'ThetaToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D50').Value = '2.59'
$ws.Range('E50').Value = '  -9.41%  '
$ws.Range('D51').Value = '1.00'
$ws.Range('E51').Value = '  -0.06%  '
